$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.166.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "'2.052.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'248.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'0.666"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "'58.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.49%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.385"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'0.0785"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("D11").Value = "'0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "'15.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "'2.350.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "'0.833"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "'5.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.30%  "
$ws.Range("D16").Value = "'2.054.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "'18.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +17.28%  "
$ws.Range("D18").Value = "'37.164.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "'75.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'0.0₃0900"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("D21").Value = "'5.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").Value = "'237.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'2.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("D25").Value = "'2.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "'169.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "'9.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "'20.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "'0.125"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "'4.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").Value = "'1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.41%  "
$ws.Range("D32").Value = "'0.0623"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").Value = "'4.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").Value = "'0.0900"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "'2.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -4.14%  "
$ws.Range("D40").Value = "'3.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.41%  "
$ws.Range("D41").Value = "'5.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.04%  "
$ws.Range("D42").Value = "'0.0223"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "'17.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("D44").Value = "'1.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").Value = "'96.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").Value = "'2.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "'2.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'1.283.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("D49").Value = "'6.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").Value = "'2.243.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "'3.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -19.09%  "
